$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46002

$ws.Range("B2").Value = 103.52
$ws.Range("C2").Value = 93.81
$ws.Range("D2").Value = 92.33
$ws.Range("E2").Value = 82.84999999999999
$ws.Range("F2").Value = 80.45999999999999
$ws.Range("G2").Value = 82.27
$ws.Range("H2").Value = 95.09
$ws.Range("I2").Value = 113
$ws.Range("J2").Value = 118.51
$ws.Range("K2").Value = 105.23
$ws.Range("L2").Value = 85.5
$ws.Range("M2").Value = 78.77
$ws.Range("N2").Value = 76.55
$ws.Range("O2").Value = 75.02
$ws.Range("P2").Value = 75.15000000000001
$ws.Range("Q2").Value = 77.84
$ws.Range("R2").Value = 92.76000000000001
$ws.Range("S2").Value = 112.78
$ws.Range("T2").Value = 111
$ws.Range("U2").Value = 115.58
$ws.Range("V2").Value = 124.52
$ws.Range("W2").Value = 118
$ws.Range("X2").Value = 103.29
$ws.Range("Y2").Value = 94.75
$ws.Range("Z2").Value = 96.19

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 110.14

$ws.Range("AD2").Value = 121.26

$ws.Range("AF2").Value = 113.29

$ws.Range("AG2").Value = "1h-23h"
